$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 35.73885133333334
$ws.Cells.Item(2, 8).Value = 107.216554
$ws.Cells.Item(2, 9).Value = 0.01949729408921566
$ws.Cells.Item(2, 10).Value = 0.01949729408921566
$ws.Cells.Item(2, 13).Value = 11.128273
$ws.Cells.Item(2, 14).Value = 33.384819
$ws.Cells.Item(2, 15).Value = 0.1975004092010595
$ws.Cells.Item(2, 16).Value = 0.1975004092010595
$ws.Cells.Item(2, 17).Value = 397.7116943437474
$ws.Cells.Item(2, 18).Value = 3579.405249093726
$ws.Cells.Item(2, 19).Value = 0.003850723560933491
$ws.Cells.Item(2, 20).Value = 0.003850723560933491
$ws.Cells.Item(3, 7).Value = 35.73885133333334
$ws.Cells.Item(3, 8).Value = 107.216554
$ws.Cells.Item(3, 9).Value = 0.01949729408921566
$ws.Cells.Item(3, 10).Value = 0.01949729408921566
$ws.Cells.Item(3, 15).Value = 0.1348505024127659
$ws.Cells.Item(3, 16).Value = 0.1348505024127659
$ws.Cells.Item(3, 17).Value = 271.5519528017211
$ws.Cells.Item(3, 18).Value = 2443.96757521549
$ws.Cells.Item(3, 19).Value = 0.002629219903620182
$ws.Cells.Item(3, 20).Value = 0.002629219903620182
$ws.Cells.Item(4, 7).Value = 35.73885133333334
$ws.Cells.Item(4, 8).Value = 107.216554
$ws.Cells.Item(4, 9).Value = 0.01949729408921566
$ws.Cells.Item(4, 10).Value = 0.01949729408921566
$ws.Cells.Item(4, 13).Value = 8.642352000000001
$ws.Cells.Item(4, 14).Value = 25.927056
$ws.Cells.Item(4, 15).Value = 0.153381217054937
$ws.Cells.Item(4, 16).Value = 0.153381217054937
$ws.Cells.Item(4, 17).Value = 308.8677332983361
$ws.Cells.Item(4, 18).Value = 2779.809599685024
$ws.Cells.Item(4, 19).Value = 0.002990518696681927
$ws.Cells.Item(4, 20).Value = 0.002990518696681926
$ws.Cells.Item(5, 7).Value = 35.73885133333334
$ws.Cells.Item(5, 8).Value = 107.216554
$ws.Cells.Item(5, 9).Value = 0.01949729408921566
$ws.Cells.Item(5, 10).Value = 0.01949729408921566
$ws.Cells.Item(5, 13).Value = 28.976716
$ws.Cells.Item(5, 14).Value = 86.930148
$ws.Cells.Item(5, 15).Value = 0.5142678713312377
$ws.Cells.Item(5, 16).Value = 0.5142678713312377
$ws.Cells.Item(5, 17).Value = 1035.594545252221
$ws.Cells.Item(5, 18).Value = 9320.350907269993
$ws.Cells.Item(5, 19).Value = 0.01002683192798006
$ws.Cells.Item(5, 20).Value = 0.01002683192798006
$ws.Cells.Item(6, 8).Value = 5067.86792
$ws.Cells.Item(6, 9).Value = 0.9215900675332435
$ws.Cells.Item(6, 10).Value = 0.9215900675332435
$ws.Cells.Item(6, 13).Value = 11.128273
$ws.Cells.Item(6, 14).Value = 33.384819
$ws.Cells.Item(6, 15).Value = 0.1975004092010595
$ws.Cells.Item(6, 16).Value = 0.1975004092010595
$ws.Cells.Item(6, 17).Value = 18798.87258056739
$ws.Cells.Item(6, 18).Value = 169189.8532251065
$ws.Cells.Item(6, 19).Value = 0.1820144154534476
$ws.Cells.Item(6, 20).Value = 0.1820144154534476
$ws.Cells.Item(7, 8).Value = 5067.86792
$ws.Cells.Item(7, 9).Value = 0.9215900675332435
$ws.Cells.Item(7, 10).Value = 0.9215900675332435
$ws.Cells.Item(7, 15).Value = 0.1348505024127659
$ws.Cells.Item(7, 16).Value = 0.1348505024127659
$ws.Cells.Item(7, 19).Value = 0.1242768836254727
$ws.Cells.Item(7, 20).Value = 0.1242768836254727
$ws.Cells.Item(8, 8).Value = 5067.86792
$ws.Cells.Item(8, 9).Value = 0.9215900675332435
$ws.Cells.Item(8, 10).Value = 0.9215900675332435
$ws.Cells.Item(8, 13).Value = 8.642352000000001
$ws.Cells.Item(8, 14).Value = 25.927056
$ws.Cells.Item(8, 15).Value = 0.153381217054937
$ws.Cells.Item(8, 16).Value = 0.153381217054937
$ws.Cells.Item(8, 17).Value = 14599.43281804928
$ws.Cells.Item(8, 18).Value = 131394.8953624435
$ws.Cells.Item(8, 19).Value = 0.1413546061839905
$ws.Cells.Item(8, 20).Value = 0.1413546061839904
$ws.Cells.Item(9, 8).Value = 5067.86792
$ws.Cells.Item(9, 9).Value = 0.9215900675332435
$ws.Cells.Item(9, 10).Value = 0.9215900675332435
$ws.Cells.Item(9, 13).Value = 28.976716
$ws.Cells.Item(9, 14).Value = 86.930148
$ws.Cells.Item(9, 15).Value = 0.5142678713312377
$ws.Cells.Item(9, 16).Value = 0.5142678713312377
$ws.Cells.Item(9, 17).Value = 48950.0564811169
$ws.Cells.Item(9, 18).Value = 440550.5083300521
$ws.Cells.Item(9, 19).Value = 0.4739441622703328
$ws.Cells.Item(9, 20).Value = 0.4739441622703328
$ws.Cells.Item(10, 7).Value = 93.641553
$ws.Cells.Item(10, 8).Value = 280.924659
$ws.Cells.Item(10, 9).Value = 0.05108605424341119
$ws.Cells.Item(10, 10).Value = 0.05108605424341119
$ws.Cells.Item(10, 13).Value = 11.128273
$ws.Cells.Item(10, 14).Value = 33.384819
$ws.Cells.Item(10, 15).Value = 0.1975004092010595
$ws.Cells.Item(10, 16).Value = 0.1975004092010595
$ws.Cells.Item(10, 17).Value = 1042.068765927969
$ws.Cells.Item(10, 18).Value = 9378.618893351722
$ws.Cells.Item(10, 19).Value = 0.01008951661754123
$ws.Cells.Item(10, 20).Value = 0.01008951661754123
$ws.Cells.Item(11, 7).Value = 93.641553
$ws.Cells.Item(11, 8).Value = 280.924659
$ws.Cells.Item(11, 9).Value = 0.05108605424341119
$ws.Cells.Item(11, 10).Value = 0.05108605424341119
$ws.Cells.Item(11, 15).Value = 0.1348505024127659
$ws.Cells.Item(11, 16).Value = 0.1348505024127659
$ws.Cells.Item(11, 17).Value = 711.5099011819351
$ws.Cells.Item(11, 18).Value = 6403.589110637416
$ws.Cells.Item(11, 19).Value = 0.006888980081009809
$ws.Cells.Item(11, 20).Value = 0.006888980081009809
$ws.Cells.Item(12, 7).Value = 93.641553
$ws.Cells.Item(12, 8).Value = 280.924659
$ws.Cells.Item(12, 9).Value = 0.05108605424341119
$ws.Cells.Item(12, 10).Value = 0.05108605424341119
$ws.Cells.Item(12, 13).Value = 8.642352000000001
$ws.Cells.Item(12, 14).Value = 25.927056
$ws.Cells.Item(12, 15).Value = 0.153381217054937
$ws.Cells.Item(12, 16).Value = 0.153381217054937
$ws.Cells.Item(12, 17).Value = 809.2832628526561
$ws.Cells.Item(12, 18).Value = 7283.549365673905
$ws.Cells.Item(12, 19).Value = 0.007835641174388936
$ws.Cells.Item(12, 20).Value = 0.007835641174388935
$ws.Cells.Item(13, 7).Value = 93.641553
$ws.Cells.Item(13, 8).Value = 280.924659
$ws.Cells.Item(13, 9).Value = 0.05108605424341119
$ws.Cells.Item(13, 10).Value = 0.05108605424341119
$ws.Cells.Item(13, 13).Value = 28.976716
$ws.Cells.Item(13, 14).Value = 86.930148
$ws.Cells.Item(13, 15).Value = 0.5142678713312377
$ws.Cells.Item(13, 16).Value = 0.5142678713312377
$ws.Cells.Item(13, 17).Value = 2713.424687079948
$ws.Cells.Item(13, 18).Value = 24420.82218371953
$ws.Cells.Item(13, 19).Value = 0.02627191637047122
$ws.Cells.Item(13, 20).Value = 0.02627191637047122
$ws.Cells.Item(14, 7).Value = 14.34625366666667
$ws.Cells.Item(14, 8).Value = 43.038761
$ws.Cells.Item(14, 9).Value = 0.007826584134129748
$ws.Cells.Item(14, 10).Value = 0.007826584134129748
$ws.Cells.Item(14, 13).Value = 11.128273
$ws.Cells.Item(14, 14).Value = 33.384819
$ws.Cells.Item(14, 15).Value = 0.1975004092010595
$ws.Cells.Item(14, 16).Value = 0.1975004092010595
$ws.Cells.Item(14, 17).Value = 159.6490273299177
$ws.Cells.Item(14, 18).Value = 1436.841245969259
$ws.Cells.Item(14, 19).Value = 0.001545753569137145
$ws.Cells.Item(14, 20).Value = 0.001545753569137145
$ws.Cells.Item(15, 7).Value = 14.34625366666667
$ws.Cells.Item(15, 8).Value = 43.038761
$ws.Cells.Item(15, 9).Value = 0.007826584134129748
$ws.Cells.Item(15, 10).Value = 0.007826584134129748
$ws.Cells.Item(15, 15).Value = 0.1348505024127659
$ws.Cells.Item(15, 16).Value = 0.1348505024127659
$ws.Cells.Item(15, 17).Value = 109.0061110872539
$ws.Cells.Item(15, 18).Value = 981.0549997852851
$ws.Cells.Item(15, 19).Value = 0.001055418802663179
$ws.Cells.Item(15, 20).Value = 0.001055418802663179
$ws.Cells.Item(16, 7).Value = 14.34625366666667
$ws.Cells.Item(16, 8).Value = 43.038761
$ws.Cells.Item(16, 9).Value = 0.007826584134129748
$ws.Cells.Item(16, 10).Value = 0.007826584134129748
$ws.Cells.Item(16, 13).Value = 8.642352000000001
$ws.Cells.Item(16, 14).Value = 25.927056
$ws.Cells.Item(16, 15).Value = 0.153381217054937
$ws.Cells.Item(16, 16).Value = 0.153381217054937
$ws.Cells.Item(16, 17).Value = 123.985374068624
$ws.Cells.Item(16, 18).Value = 1115.868366617616
$ws.Cells.Item(16, 19).Value = 0.001200450999875681
$ws.Cells.Item(16, 20).Value = 0.001200450999875681
$ws.Cells.Item(17, 7).Value = 14.34625366666667
$ws.Cells.Item(17, 8).Value = 43.038761
$ws.Cells.Item(17, 9).Value = 0.007826584134129748
$ws.Cells.Item(17, 10).Value = 0.007826584134129748
$ws.Cells.Item(17, 13).Value = 28.976716
$ws.Cells.Item(17, 14).Value = 86.930148
$ws.Cells.Item(17, 15).Value = 0.5142678713312377
$ws.Cells.Item(17, 16).Value = 0.5142678713312377
$ws.Cells.Item(17, 17).Value = 415.7073181629587
$ws.Cells.Item(17, 18).Value = 3741.365863466628
$ws.Cells.Item(17, 19).Value = 0.004024960762453744
$ws.Cells.Item(17, 20).Value = 0.004024960762453744
